$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F3: plain numeric cell (count), stays a number
$ws.Range("F3").Value = 58

# F4 and F11:F63: "Laenge" values are stored as TEXT (shared strings) in the
# original workbook, not numbers. Pre-format the range as Text ("@") before
# writing so Excel does not silently coerce the numeric-looking strings into
# real numbers, then clear the formatting again so no extra cell style is left
# behind (matches the un-styled cells in the source file).
$textCells = $ws.Range("F4,F11,F12,F13,F14,F15,F16,F17,F18,F19,F20,F22,F23,F24,F27,F28,F29,F30,F31,F33,F34,F35,F36,F37,F38,F40,F43,F44,F45,F47,F48,F49,F50,F51,F52,F53,F54,F55,F58,F60,F61,F62,F63")
$textCells.NumberFormat = "@"

$ws.Range("F4").Value = "80.13"
$ws.Range("F11").Value = "13.85"
$ws.Range("F12").Value = "14.31"
$ws.Range("F13").Value = "14.55"
$ws.Range("F14").Value = "15.02"
$ws.Range("F15").Value = "14.78"
$ws.Range("F16").Value = "6.55"
$ws.Range("F17").Value = "15.25"
$ws.Range("F18").Value = "13.42"
$ws.Range("F19").Value = "6.6"
$ws.Range("F20").Value = "6.96"
$ws.Range("F22").Value = "6.6"
$ws.Range("F23").Value = "15.7"
$ws.Range("F24").Value = "6.11"
$ws.Range("F27").Value = "6.92"
$ws.Range("F28").Value = "16.32"
$ws.Range("F29").Value = "6.11"
$ws.Range("F30").Value = "6.21"
$ws.Range("F31").Value = "13.9"
$ws.Range("F33").Value = "6.21"
$ws.Range("F34").Value = "14.05"
$ws.Range("F35").Value = "6.8"
$ws.Range("F36").Value = "6.8"
$ws.Range("F37").Value = "18.21"
$ws.Range("F38").Value = "14.1"
$ws.Range("F40").Value = "7.57"
$ws.Range("F43").Value = "8.21"
$ws.Range("F44").Value = "8.21"
$ws.Range("F45").Value = "21.48"
$ws.Range("F47").Value = "6.91"
$ws.Range("F48").Value = "13.59"
$ws.Range("F49").Value = "6.88"
$ws.Range("F50").Value = "7.37"
$ws.Range("F51").Value = "7.37"
$ws.Range("F52").Value = "6.88"
$ws.Range("F53").Value = "16.41"
$ws.Range("F54").Value = "13.81"
$ws.Range("F55").Value = "6.97"
$ws.Range("F58").Value = "13.96"
$ws.Range("F60").Value = "7.5"
$ws.Range("F61").Value = "18.3"
$ws.Range("F62").Value = "7.78"
$ws.Range("F63").Value = "7.78"

$textCells.ClearFormats()
